# Updates cryptos list prices/volumes (mirrors "Updated cryptos list ... with GitHub Actions").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (Price text for column D, Volume(1h) text for column E).
# $null means "leave that column untouched" for that row.
$updates = @{
    2  = @("60.565.17", "  +0.13%  ")
    3  = @("2.595.05",  "  +0.14%  ")
    4  = @($null,       "  +0.00%  ")
    5  = @("516.94",    "  +1.90%  ")
    6  = @("153.63",    "  -0.22%  ")
    7  = @($null,       "  +0.04%  ")
    8  = @($null,       "  +2.99%  ")
    9  = @("6.70",      "  +0.91%  ")
    10 = @($null,       "  +1.69%  ")
    11 = @("0.346",     "  +0.10%  ")
    13 = @("3.049.89",  "  +0.18%  ")
    14 = @("60.581.20", "  +0.20%  ")
    15 = @("21.66",     "  +0.05%  ")
    17 = @("2.600.97",  "  +0.37%  ")
    18 = @($null,       "  -1.66%  ")
    19 = @("350.80",    "  +0.96%  ")
    20 = @("10.55",     "  +1.73%  ")
    21 = @("6.21",      "  +1.31%  ")
    22 = @($null,       "  -0.02%  ")
    23 = @("60.95",     "  +1.09%  ")
    24 = @($null,       "  +1.37%  ")
    25 = @($null,       "  -0.07%  ")
    26 = @("2.712.95",  $null)
    27 = @("0.999",     "  +0.33%  ")
    28 = @("0.0₃0842",  "  -0.66%  ")
    29 = @($null,       "  -1.68%  ")
    30 = @($null,       "  +0.03%  ")
    31 = @($null,       "  +8.77%  ")
    32 = @($null,       "  +0.12%  ")
    33 = @($null,       "  +2.40%  ")
    34 = @("149.67",    "  -3.11%  ")
    35 = @("4.14",      "  +3.23%  ")
    36 = @($null,       "  +0.33%  ")
    37 = @("0.917",     "  +8.33%  ")
    38 = @($null,       "  +1.53%  ")
    39 = @($null,       "  +0.32%  ")
    40 = @("36.34",     "  +1.34%  ")
    41 = @("0.837",     "  -1.55%  ")
    42 = @("286.45",    "  -3.69%  ")
    43 = @($null,       "  +1.44%  ")
    44 = @($null,       "  +0.49%  ")
    45 = @("0.0559",    "  -0.59%  ")
    46 = @("0.997",     "  +0.02%  ")
    47 = @("19.52",     "  -0.80%  ")
    48 = @($null,       "  +0.87%  ")
    49 = @($null,       "  -2.21%  ")
    50 = @($null,       "  +0.13%  ")
    51 = @("18.97",     "  +6.95%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $priceVal = $vals[0]
    $volVal = $vals[1]

    if ($priceVal -ne $null) {
        # Leading apostrophe forces the numeric-looking text to stay text
        # (same as typing '60.565.17 into Excel) instead of being parsed
        # into a number.
        $ws.Cells.Item($row, 4).Value = "'" + $priceVal
    }
    if ($volVal -ne $null) {
        $ws.Cells.Item($row, 5).Value = $volVal
    }
}
